$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 (NonFastTag) edits ---

# Update header labels (G1/H1): "Policy Start Date"/"Policy End Date" -> "Policy Start"/"Policy End"
$ws2.Range("G1").Value = "Policy Start"
$ws2.Range("H1").Value = "Policy End"

# Insert a new row above the existing row 2 (shifts current rows 2-3 down to 3-4)
$ws2.Cells.Item(2, 1).EntireRow.Insert()

# New row 2: TVS / Jupyter 110 policy, now tagged with a different registration number
$ws2.Range("A2").Value = "KL38P5008"
$ws2.Range("B2").Value = "TVS"
$ws2.Range("C2").Value = 2016
$ws2.Range("D2").Value = "Jupyter 110"
$ws2.Range("E2").Value = 560
$ws2.Range("F2").Value = "BADAA"
$ws2.Range("G2").Value = 44854
$ws2.Range("G2").NumberFormat = $ws2.Range("G3").NumberFormat()
$ws2.Range("H2").Value = 45219
$ws2.Range("H2").NumberFormat = $ws2.Range("H3").NumberFormat()
$ws2.Range("I2").Value = "Third Party"
$ws2.Range("K2").Value = "Not Expired"

# Row 3 keeps the original Hero/Xpulse data (was row 2 before the insert) - unchanged

# Row 4 currently holds the old KA19P8488/TVS/Jupyter duplicate row pushed down by the
# insert; replace it with the new Bajaj/Pulser policy record
$ws2.Range("A4").Value = "KL01AU585"
$ws2.Range("B4").Value = "Bajaj"
$ws2.Range("C4").Value = 2013
$ws2.Range("D4").Value = "Pulser"
$ws2.Range("E4").Value = 850
$ws2.Range("F4").Value = "AAA"
$ws2.Range("G4").Value = 42849
$ws2.Range("G4").NumberFormat = $ws2.Range("G3").NumberFormat()
$ws2.Range("H4").Value = 43214
$ws2.Range("H4").NumberFormat = $ws2.Range("H3").NumberFormat()
$ws2.Range("I4").Value = "Comprehensive"
$ws2.Range("K4").Value = "Expired"

# Row 5 (new): LML / Freedom policy
$ws2.Range("A5").Value = "GJ05KP2603"
$ws2.Range("B5").Value = "LML"
$ws2.Range("C5").Value = 2010
$ws2.Range("D5").Value = "Freedom"
$ws2.Range("E5").Value = 730
$ws2.Range("F5").Value = "BADAA"
$ws2.Range("G5").Value = 44854
$ws2.Range("G5").NumberFormat = $ws2.Range("G3").NumberFormat()
$ws2.Range("H5").Value = 45280
$ws2.Range("H5").NumberFormat = $ws2.Range("H3").NumberFormat()
$ws2.Range("I5").Value = "Third Party"
$ws2.Range("K5").Value = "Not Expired"

# Make NonFastTag the active sheet/tab, with the selection on N11
$ws2.Activate()
$ws2.Range("N11").Select()

# Sheet1 (FastTag) keeps its own selection (C23) but is no longer the active/selected tab;
# activating ws2 above already moved tabSelected off of sheet1's sheetView.
